$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Drum" (DRUMMOND?) sowing-date / treatment labels to "Stam"
# (STAMPEDE) for the 4/6-leaf photoperiod phenology rework.
# Column I (rows 11-19) holds the long SimulationName values,
# column L (rows 11-19) holds the short treatment code used in lookups.

$ws.Range("I11").Value = "SowingDate2002SowStamMar"
$ws.Range("L11").Value = "StamMar"

$ws.Range("I12").Value = "SowingDate2002SowStamMay"
$ws.Range("L12").Value = "StamMay"

$ws.Range("I13").Value = "SowingDate2002SowStamSep"
$ws.Range("L13").Value = "StamSep"

$ws.Range("I14").Value = "SowingDate2003SowStamMar"
$ws.Range("L14").Value = "StamMar"

$ws.Range("I15").Value = "SowingDate2003SowStamMay"
$ws.Range("L15").Value = "StamMay"

$ws.Range("I16").Value = "SowingDate2003SowStamSep"
$ws.Range("L16").Value = "StamSep"

$ws.Range("I17").Value = "SowingDate2004SowStamMar"
$ws.Range("L17").Value = "StamMar"

$ws.Range("I18").Value = "SowingDate2004SowStamMay"
$ws.Range("L18").Value = "StamMay"

$ws.Range("I19").Value = "SowingDate2004SowStamSep"
$ws.Range("L19").Value = "StamSep"

# Select the whole sheet (mirrors the "select all" state captured in the
# saved view) instead of the previous E9:E10 selection.
$ws.Range("A1:XFD1048576").Select()
